$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2-10), following the updated NATMI computation
# (Ligand-/Receptor-expressing cell counts changed from 1 to 3, with
# recalculated average/total expression values and specificities).

$data = @(
    @{ Row=2;  E=3; G=19.952791;           H=59.858373;          I=0.6991080337323071; J=0.6991080337323072; K=3; M=2.300909333333333;  N=6.902728;   O=0.03776979643482627; P=0.03776979643482627; Q=45.90956303794933;  R=413.186067341544;   S=0.0264051681200209;   T=0.0264051681200209 }
    @{ Row=3;  E=3; G=19.952791;           H=59.858373;          I=0.6991080337323071; J=0.6991080337323072; K=3; M=16.96312166666667;  N=50.889365;  O=0.2784523679257784;  P=0.2784523679257784;  Q=338.4616213225717;  R=3046.154591903145;  S=0.1946682874286959;   T=0.1946682874286959 }
    @{ Row=4;  E=3; G=19.952791;           H=59.858373;          I=0.6991080337323071; J=0.6991080337323072; K=3; M=41.655263;          N=124.965789; O=0.6837778356393953;  P=0.6837778356393953;  Q=831.138756689033;   R=7480.248810201297;  S=0.4780345781835904;   T=0.4780345781835904 }
    @{ Row=5;  E=3; G=3.140511;            H=9.421533;           I=0.1100375616686749; J=0.1100375616686749; K=3; M=2.300909333333333;  N=6.902728;   O=0.03776979643482627; P=0.03776979643482627; Q=7.226031071336;     R=65.034279642024;    S=0.004156096304410493; T=0.004156096304410494 }
    @{ Row=6;  E=3; G=3.140511;            H=9.421533;           I=0.1100375616686749; J=0.1100375616686749; K=3; M=16.96312166666667;  N=50.889365;  O=0.2784523679257784;  P=0.2784523679257784;  Q=53.272870188505;    R=479.455831696545;   S=0.0306402196074214;   T=0.0306402196074214 }
    @{ Row=7;  E=3; G=3.140511;            H=9.421533;           I=0.1100375616686749; J=0.1100375616686749; K=3; M=41.655263;          N=124.965789; O=0.6837778356393953;  P=0.6837778356393953;  Q=130.818811659393;   R=1177.369304934537;  S=0.07524124575684302;  T=0.07524124575684303 }
    @{ Row=8;  E=3; G=5.447052333333334;   H=16.341157;          I=0.1908544045990179; J=0.1908544045990179; K=3; M=2.300909333333333;  N=6.902728;   O=0.03776979643482627; P=0.03776979643482627; Q=12.53317355292178;  R=112.798561976296;   S=0.007208532010394877; T=0.007208532010394877 }
    @{ Row=9;  E=3; G=5.447052333333334;   H=16.341157;          I=0.1908544045990179; J=0.1908544045990179; K=3; M=16.96312166666667;  N=50.889365;  O=0.2784523679257784;  P=0.2784523679257784;  Q=92.3990114550339;   R=831.5911030953051;  S=0.05314386088966111;  T=0.05314386088966112 }
    @{ Row=10; E=3; G=5.447052333333334;   H=16.341157;          I=0.1908544045990179; J=0.1908544045990179; K=3; M=41.655263;          N=124.965789; O=0.6837778356393953;  P=0.6837778356393953;  Q=226.8983975197637;  R=2042.085577677873;  S=0.1305020116989619;   T=0.1305020116989619 }
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Range("E$r").Value = $rowData.E
    $ws.Range("G$r").Value = $rowData.G
    $ws.Range("H$r").Value = $rowData.H
    $ws.Range("I$r").Value = $rowData.I
    $ws.Range("J$r").Value = $rowData.J
    $ws.Range("K$r").Value = $rowData.K
    $ws.Range("M$r").Value = $rowData.M
    $ws.Range("N$r").Value = $rowData.N
    $ws.Range("O$r").Value = $rowData.O
    $ws.Range("P$r").Value = $rowData.P
    $ws.Range("Q$r").Value = $rowData.Q
    $ws.Range("R$r").Value = $rowData.R
    $ws.Range("S$r").Value = $rowData.S
    $ws.Range("T$r").Value = $rowData.T
}
